# Update the "Förändrad" (Changed) date column (C) for rows 2-21:
# each date serial value increases by one day (45184 -> 45185),
# i.e. 2023-09-15 -> 2023-09-16. Formatting/style is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 21; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45185
    }
}
